# Horarios actualizados Linea 141 - 53
# Refresh the scraped schedule tables on all three sheets with the new
# scrape snapshot (timestamp 04:01:06), which grows the main tables from
# 12 to 18 data rows (and the filtered "215" sheet from 5 to 6 rows).

$wb = $excel.ActiveWorkbook

$lastUpdate = "Última actualización: 04:01:06"

# Full refreshed dataset (sorted by Minutos ascending) shared by the
# "LP1912" and "6203-6173" sheets.
$mainData = @(
    @("04:01:06","04:02","81_EL PELIGRO",1),
    @("00:46:06","01:12","215_ALUAR",26),
    @("04:01:06","04:47","215_EL PELIGRO",46),
    @("03:46:12","04:46","215A_EL PATO",60),
    @("01:55:38","03:02","15_ABASTO",67),
    @("04:01:06","05:12","17_ROMERO",71),
    @("00:46:06","01:58","14_ABASTO",72),
    @("03:46:12","05:16","17_ROMERO",90),
    @("04:01:06","05:32","81_EL PELIGRO",91),
    @("02:29:13","04:01","81_EL PELIGRO",92),
    @("03:46:12","05:22","23_HERNANDEZ",96),
    @("01:22:42","02:58","215_ALUAR",96),
    @("04:01:06","05:45","14_ABASTO",104),
    @("03:46:12","05:35","215B_EL PATO",109),
    @("04:01:06","05:52","17_ROMERO",111),
    @("01:55:38","03:48","14_ABASTO",113),
    @("03:00:53","04:53","11_ETCHEVERRY",113),
    @("02:47:42","04:45","215A_EL PATO",118)
)

# Refreshed dataset for the "LP1912-215" sheet (subset of the main data
# whose "Linea" contains "215").
$sheet215Data = @(
    @("00:46:06","01:12","215_ALUAR",26),
    @("04:01:06","04:47","215_EL PELIGRO",46),
    @("03:46:12","04:46","215A_EL PATO",60),
    @("01:22:42","02:58","215_ALUAR",96),
    @("03:46:12","05:35","215B_EL PATO",109),
    @("02:47:42","04:45","215A_EL PATO",118)
)

function Write-Table($ws, $rows) {
    $ws.Range("A2").Value = $lastUpdate
    $ws.Cells.Item(3, 1).Value = "Total filas: " + $rows.Count

    $r = 6
    foreach ($row in $rows) {
        $ws.Cells.Item($r, 1).Value = $row[0]
        $ws.Cells.Item($r, 2).Value = $row[1]
        $ws.Cells.Item($r, 3).Value = $row[2]
        $ws.Cells.Item($r, 4).Value = $row[3]
        $r = $r + 1
    }
}

$wsMain = $wb.Worksheets.Item("LP1912")
Write-Table $wsMain $mainData

$ws215 = $wb.Worksheets.Item("LP1912-215")
Write-Table $ws215 $sheet215Data

$wsOther = $wb.Worksheets.Item("6203-6173")
Write-Table $wsOther $mainData
